# Refresh the scraped crypto price/volume snapshot (GitHub Actions bot run).
# Column D ("Price") and column E ("Volume(1h)") are plain text cells (the
# source feed renders pre-formatted strings, e.g. "64.321.20" / "  -4.50%  "),
# so every write below targets the Coin/Link/Price/Volume columns as text.
$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Row 2: Bitcoin
$ws.Range('D2').Value = '64.321.20'
$ws.Range('E2').Value = '  -4.50%  '

# Row 3: Ethereum
$ws.Range('D3').Value = '3.412.25'
$ws.Range('E3').Value = '  -3.87%  '

# Row 5: BNB
$ws.Range('D5').Value = '''560.97'
$ws.Range('E5').Value = '  -0.04%  '

# Row 6: Solana
$ws.Range('D6').Value = '''173.16'
$ws.Range('E6').Value = '  -9.00%  '

# Row 7: XRP
$ws.Range('D7').Value = '''0.620'
$ws.Range('E7').Value = '  +0.48%  '

# Row 8: USDC
$ws.Range('E8').Value = '  +0.02%  '

# Row 9: Cardano
$ws.Range('D9').Value = '''0.619'
$ws.Range('E9').Value = '  -1.93%  '

# Row 10: Dogecoin
$ws.Range('D10').Value = '''0.154'
$ws.Range('E10').Value = '  +0.41%  '

# Row 11: Avalanche
$ws.Range('D11').Value = '''56.15'
$ws.Range('E11').Value = '  +1.91%  '

# Row 12: ShibaInu
$ws.Range('D12').Value = '''0.0000269'
$ws.Range('E12').Value = '  -1.18%  '

# Row 13: Polkadot
$ws.Range('D13').Value = '''9.03'
$ws.Range('E13').Value = '  -3.71%  '

# Row 14: WrappedliquidstakedEther2.0
$ws.Range('D14').Value = '3.956.21'
$ws.Range('E14').Value = '  -3.62%  '

# Row 15: TRON
$ws.Range('E15').Value = '  -1.25%  '

# Row 16: WrappedEther
$ws.Range('D16').Value = '3.410.65'
$ws.Range('E16').Value = '  -3.74%  '

# Row 17: Chainlink
$ws.Range('D17').Value = '''17.98'
$ws.Range('E17').Value = '  -1.65%  '

# Row 18: Uniswap
$ws.Range('D18').Value = '''11.77'
$ws.Range('E18').Value = '  -2.40%  '

# Row 19: WrappedBTC
$ws.Range('D19').Value = '64.296.31'
$ws.Range('E19').Value = '  -4.42%  '

# Row 20: Polygon
$ws.Range('D20').Value = '''0.988'
$ws.Range('E20').Value = '  -1.24%  '

# Row 21: BitcoinCash
$ws.Range('D21').Value = '''407.93'
$ws.Range('E21').Value = '  -5.12%  '

# Row 22: PancakeSwap
$ws.Range('D22').Value = '''4.12'
$ws.Range('E22').Value = '  +0.24%  '

# Row 23: Toncoin
$ws.Range('D23').Value = '''4.38'
$ws.Range('E23').Value = '  +5.53%  '

# Row 24: InternetComputer(DFINITY)
$ws.Range('D24').Value = '''13.33'
$ws.Range('E24').Value = '  +7.80%  '

# Row 25: Litecoin
$ws.Range('D25').Value = '''82.88'
$ws.Range('E25').Value = '  -2.61%  '

# Row 26: RenderToken
$ws.Range('D26').Value = '''10.71'
$ws.Range('E26').Value = '  -3.14%  '

# Row 27: ImmutableX
$ws.Range('E27').Value = '  -5.17%  '

# Row 28: Filecoin
$ws.Range('D28').Value = '''8.81'
$ws.Range('E28').Value = '  -2.51%  '

# Row 29: EthereumClassic
$ws.Range('D29').Value = '''29.57'
$ws.Range('E29').Value = '  -3.19%  '

# Row 30: NEARProtocol
$ws.Range('D30').Value = '''6.60'
$ws.Range('E30').Value = '  -0.56%  '

# Row 31: Bittensor
$ws.Range('D31').Value = '''591.33'
$ws.Range('E31').Value = '  -7.61%  '

# Row 32: Cosmos
$ws.Range('D32').Value = '''11.44'
$ws.Range('E32').Value = '  -2.52%  '

# Row 33: Hedera
$ws.Range('E33').Value = '  -4.00%  '

# Row 34: Kaspa
$ws.Range('B34').Value = 'OKB'
$ws.Range('C34').Value = 'https://coinranking.com/coin/PDKcptVnzJTmN+okb-okb'
$ws.Range('D34').Value = '''58.87'
$ws.Range('E34').Value = '  -2.18%  '

# Row 35: OKB
$ws.Range('B35').Value = 'Kaspa'
$ws.Range('C35').Value = 'https://coinranking.com/coin/V8GxkwWow+kaspa-kas'
$ws.Range('D35').Value = '''0.153'
$ws.Range('E35').Value = '  +3.55%  '

# Row 36: Dai
$ws.Range('E36').Value = '  +0.18%  '

# Row 37: InjectiveProtocol
$ws.Range('D37').Value = '''35.72'
$ws.Range('E37').Value = '  -7.24%  '

# Row 38: Stacks
$ws.Range('D38').Value = '''3.40'
$ws.Range('E38').Value = '  -0.30%  '

# Row 39: TheGraph
$ws.Range('D39').Value = '''0.371'
$ws.Range('E39').Value = '  -4.60%  '

# Row 40: PEPE
$ws.Range('D40').Value = '0.0₃0735'
$ws.Range('E40').Value = '  -9.74%  '

# Row 41: Maker
$ws.Range('D41').Value = '3.179.42'
$ws.Range('E41').Value = '  +2.33%  '

# Row 42: FirstDigitalUSD
$ws.Range('E42').Value = '  +0.05%  '

# Row 43: ThetaToken
$ws.Range('D43').Value = '''2.88'
$ws.Range('E43').Value = '  +0.17%  '

# Row 44: Fetch.AI
$ws.Range('D44').Value = '''2.52'
$ws.Range('E44').Value = '  -5.41%  '

# Row 45: ApeXProtocol
$ws.Range('D45').Value = '''3.22'
$ws.Range('E45').Value = '  -4.21%  '

# Row 46: VeChain
$ws.Range('D46').Value = '''0.0406'
$ws.Range('E46').Value = '  -3.20%  '

# Row 47: WEMIXToken
$ws.Range('D47').Value = '''2.62'
$ws.Range('E47').Value = '  -5.49%  '

# Row 48: Stellar
$ws.Range('E48').Value = '  -1.98%  '

# Row 49: THORChain
$ws.Range('D49').Value = '''8.26'
$ws.Range('E49').Value = '  -4.23%  '

# Row 50: Monero
$ws.Range('D50').Value = '''134.86'
$ws.Range('E50').Value = '  -4.55%  '

# Row 51: LidoDAOToken
$ws.Range('D51').Value = '''2.79'
$ws.Range('E51').Value = '  +2.34%  '
